$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(42602.009479166663, "Noun", 2723, 75, 8, 2, 1, 66, 33, 4, 3, 57, 42),
    @(42602.014409722222, "Noun", 3546, 75, 8, 2, 1, 66, 33, 4, 3, 57, 42),
    @(42602.015254629630, "Noun", 2879, 75, 8, 2, 1, 66, 33, 4, 3, 57, 42),
    @(42602.495185185187, "Noun", 3123, 75, 8, 2, 1, 66, 33, 4, 3, 57, 42),
    @(42602.495462962965, "Noun", 2934, 75, 8, 2, 1, 66, 33, 4, 3, 57, 42),
    @(42602.495717592596, "Noun", 2560, 74, 8, 1, 1, 50, 50, 2, 3, 40, 60),
    @(42602.495787037034, "Noun", 1868, 75, 8, 2, 1, 66, 33, 4, 3, 57, 42)
)

$startRow = 17
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $cell = $ws.Cells.Item($row, $c + 1)
        $cell.Value = $rowData[$c]
    }
    # Apply date number format to column A to match existing rows (style index 1 / numFmtId 22)
    $ws.Cells.Item($row, 1).NumberFormat = "m/d/yy h:mm"
}
